$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "Buy Value in GBP"
$ws.Range("G1").Value = "Sell Value in GBP"
$ws.Range("J1").Value = "Fee Value in GBP"

$headerRange = $ws.Range("A1:M1")
$headerRange.Borders.Color = 16777215
$headerRange.Borders.LineStyle = 1
